# ENW.xlsx - "Test Cases" sheet: update existing test-case row and add a new
# test-case row for the sendtoendnote scripts (ENW001 / ENW002).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Add row 3, inheriting row 2's layout/formatting (borders, wrap, etc.) ---
$ws.Range("A2:E2").Copy($ws.Range("A3:E3"))

# --- Row 2: ENW001 test case content (gets the new Jira id + description) ---
$ws.Range("A2").Value = "ENW001"
$ws.Range("B2").Value = "OPQA_1679"
$ws.Range("C2").Value = "Verify that the user is able to send the only one record at a time from article,Post,Patent view Pages"
$ws.Range("D2").Value = "Y"
$ws.Range("E2").Value = ""

# --- Row 3: new ENW002 test case ---
$ws.Range("A3").Value = "ENW002"
$ws.Range("B3").Value = "OPQA-1678"
$ws.Range("C3").Value = "Verify that user is able to send the record from below following pages: " + [char]0x2022 + "Summary lists (for Article, Patent, and Post items)"
$ws.Range("D3").Value = "Y"
$ws.Range("E3").Value = ""

# --- Row heights for the (now taller, wrapped) description rows ---
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

# --- Restore the wrap formatting on the Description column for both rows ---
$ws.Range("C2").WrapText = $true
$ws.Range("C3").WrapText = $true

# --- Cursor/selection left on C10, as in the authored workbook ---
$ws.Range("C10").Select()
